# Update cryptocurrency price/volume snapshot values in the worksheet.
# Each assignment below sets the Price (column D) and/or Volume(1h) (column E)
# text for a specific row to match the latest scraped values.
# Values that look like plain decimal numbers are prefixed with a leading
# apostrophe so Excel keeps them stored as text (matching the original
# inline-string cell type) instead of auto-converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.597.91"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.247.73"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'305.98"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'94.93"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").Value = "'34.95"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "2.589.92"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "2.240.56"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D18").Value = "44.401.62"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "0.0₃0939"
$ws.Range("E19").Value = "  -3.51%  "
$ws.Range("D20").Value = "'11.85"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").Value = "'6.20"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "'65.33"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("D28").Value = "'9.76"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "'37.26"
$ws.Range("E29").Value = "  -4.54%  "
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'19.97"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'149.81"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").Value = "'3.19"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("D39").Value = "'15.25"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("E40").Value = "  -6.76%  "
$ws.Range("D41").Value = "'3.78"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.806.91"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").Value = "'1.77"
$ws.Range("E45").Value = "  +11.32%  "
$ws.Range("D46").Value = "'81.77"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").Value = "'98.49"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("D49").Value = "'4.84"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "'68.73"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'54.03"
$ws.Range("E51").Value = "  -2.15%  "
